# Update BGDPbES BAU Guaranteed Dispatch Percentages from 0 to 1
# for: nuclear, hydro, onshore wind, solar PV, solar thermal, biomass,
# geothermal, offshore wind, and municipal solid waste.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# Rows that already contain formulas in C:AK referencing column B (directly
# or via the previous column) just need their base year (column B, 2015)
# value updated; Excel's recalculation will propagate the new value across
# the rest of the row (columns C:AK, years 2016-2050).
$ws.Range("B4").Value = 1    # nuclear
$ws.Range("B5").Value = 1    # hydro
$ws.Range("B6").Value = 1    # onshore wind
$ws.Range("B7").Value = 1    # solar PV
$ws.Range("B8").Value = 1    # solar thermal
$ws.Range("B9").Value = 1    # biomass
$ws.Range("B10").Value = 1   # geothermal
$ws.Range("B14").Value = 1   # offshore wind

# Row 17 (municipal solid waste) previously held static values with no
# formulas in C17:AK17, so set the base year value and (re)build the
# same "carry forward" formulas used by the other rows.
$ws.Range("B17").Value = 1
$ws.Range("C17").Formula = "=B17"
$ws.Range("D17:AK17").Formula = "=C17"
